$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-05 Wednesday" "2024-06-06 Thursday"

Replace-Text "364×8=2912" "204×6=1224"
Replace-Text "269×5=1345" "866×9=7794"
Replace-Text "351×2=702" "298×8=2384"
Replace-Text "459×3=1377" "333×8=2664"
Replace-Text "421×8=3368" "461×6=2766"

Replace-Text "326×2=652" "905×9=8145"
Replace-Text "711×4=2844" "110×5=550"
Replace-Text "795×3=2385" "636×8=5088"
Replace-Text "837×4=3348" "430×8=3440"
Replace-Text "718×7=5026" "505×6=3030"

Replace-Text "627×4=2508" "556×8=4448"
Replace-Text "885×3=2655" "141×7=987"
Replace-Text "676×8=5408" "623×5=3115"
Replace-Text "889×8=7112" "309×6=1854"
Replace-Text "877×4=3508" "383×7=2681"

Replace-Text "595×5=2975" "768×3=2304"
Replace-Text "511×4=2044" "312×6=1872"
Replace-Text "800×7=5600" "283×5=1415"
Replace-Text "829×3=2487" "470×6=2820"
Replace-Text "639×7=4473" "364×8=2912"

Replace-Text "958×2=1916" "390×2=780"
Replace-Text "853×9=7677" "921×2=1842"
Replace-Text "658×2=1316" "641×5=3205"
Replace-Text "329×4=1316" "435×2=870"
Replace-Text "717×8=5736" "539×6=3234"
